$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "49.564.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.57%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.646.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +6.58%  "

$ws.Range("E4").Value = "  +0.28%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "112.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.76%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "324.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.522"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.83%  "

$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.545"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.14"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.93%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.85"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.85%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0814"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.125"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.072.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.93%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.698.57"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +10.35%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.861"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.34%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "49.640.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.24%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0950"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.77%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.44%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "271.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.58%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.47"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.71%  "

$ws.Range("E27").Value = "  -0.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.12%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.58%  "

$ws.Range("E31").Value = "  -0.18%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.80"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.83%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.42"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.76%  "

$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.09%  "

$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.09"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.17%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0797"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.95"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.50%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.04"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.64%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.09"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.26%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "124.78"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.80%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.111"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.46%  "

$ws.Range("E42").Value = "  +0.34%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.98"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.63%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0315"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.87%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.099.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.89%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.25"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.26%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.78%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.92%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.56%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.27"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.34%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "58.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.91%  "
